$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030983958090183
$ws.Range("D2").Value = 1.033470337347241
$ws.Range("E2").Value = 1.030532910884497
$ws.Range("F2").Value = 1.039997854721433
$ws.Range("I2").Value = 1.033661287413533
$ws.Range("J2").Value = 1.036122189260779
$ws.Range("K2").Value = 1.036272593964153
$ws.Range("L2").Value = 1.033343647888147
$ws.Range("M2").Value = 1.042781452545055

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032530844451919
$ws.Range("D3").Value = 1.034611496746893
$ws.Range("E3").Value = 1.031869280047802
$ws.Range("F3").Value = 1.041657398511187
$ws.Range("I3").Value = 1.034075149783619
$ws.Range("J3").Value = 1.037307716553154
$ws.Range("K3").Value = 1.037221853781275
$ws.Range("L3").Value = 1.034486974139655
$ws.Range("M3").Value = 1.044249094984181

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033529646683625
$ws.Range("D4").Value = 1.035347868545169
$ws.Range("E4").Value = 1.032732261303948
$ws.Range("F4").Value = 1.042729275299207
$ws.Range("I4").Value = 1.03434033643093
$ws.Range("J4").Value = 1.038072369692494
$ws.Range("K4").Value = 1.037833483536125
$ws.Range("L4").Value = 1.035224517963617
$ws.Range("M4").Value = 1.04519634062249

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033949041696738
$ws.Range("D5").Value = 1.03565695854948
$ws.Range("E5").Value = 1.033094650124751
$ws.Range("F5").Value = 1.043179434140057
$ws.Range("I5").Value = 1.034451199188477
$ws.Range("J5").Value = 1.038393247996776
$ws.Range("K5").Value = 1.038089995106868
$ws.Range("L5").Value = 1.03553404640935
$ws.Range("M5").Value = 1.0455939938154

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034019430851429
$ws.Range("D6").Value = 1.035708828117649
$ws.Range("E6").Value = 1.033155473073587
$ws.Range("F6").Value = 1.04325499123668
$ws.Range("I6").Value = 1.034469777166203
$ws.Range("J6").Value = 1.038447090974307
$ws.Range("K6").Value = 1.038133028487197
$ws.Range("L6").Value = 1.035585986434168
$ws.Range("M6").Value = 1.045660728460562

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033535252618239
$ws.Range("D7").Value = 1.035352000501901
$ws.Range("E7").Value = 1.032737105156615
$ws.Range("F7").Value = 1.042735292130593
$ws.Range("I7").Value = 1.034341820223596
$ws.Range("J7").Value = 1.038076659561133
$ws.Range("K7").Value = 1.037836913475565
$ws.Range("L7").Value = 1.035228655991103
$ws.Range("M7").Value = 1.045201656307149

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031507184081532
$ws.Range("D8").Value = 1.033856422661779
$ws.Range("E8").Value = 1.030984907879369
$ws.Range("F8").Value = 1.040559117484385
$ws.Range("I8").Value = 1.033801696443236
$ws.Range("J8").Value = 1.036523358187135
$ws.Range("K8").Value = 1.036593943748866
$ws.Range("L8").Value = 1.033730513620601
$ws.Range("M8").Value = 1.043277955587042

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027916649987676
$ws.Range("D9").Value = 1.031205151903585
$ws.Range("E9").Value = 1.027883642375354
$ws.Range("F9").Value = 1.036708928402293
$ws.Range("I9").Value = 1.032829798743029
$ws.Range("J9").Value = 1.03376704193438
$ws.Range("K9").Value = 1.03438345083938
$ws.Range("L9").Value = 1.031072926081036
$ws.Range("M9").Value = 1.039869195297377

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025511015586107
$ws.Range("D10").Value = 1.02942656009402
$ws.Range("E10").Value = 1.025806446434116
$ws.Range("F10").Value = 1.03413103534071
$ws.Range("I10").Value = 1.032168132644771
$ws.Range("J10").Value = 1.031916106335889
$ws.Range("K10").Value = 1.032895810379662
$ws.Range("L10").Value = 1.029288865143228
$ws.Range("M10").Value = 1.037583324822663

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024466373482248
$ws.Range("D11").Value = 1.028653691150027
$ws.Range("E11").Value = 1.024904589467655
$ws.Range("F11").Value = 1.033011993070014
$ws.Range("I11").Value = 1.031878321269048
$ws.Range("J11").Value = 1.031111346425616
$ws.Range("K11").Value = 1.032248244067888
$ws.Range("L11").Value = 1.028513318877577
$ws.Range("M11").Value = 1.036590209124823

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024077885819451
$ws.Range("D12").Value = 1.028366195528604
$ws.Range("E12").Value = 1.024569226801413
$ws.Range("F12").Value = 1.032595897349332
$ws.Range("I12").Value = 1.031770171827752
$ws.Range("J12").Value = 1.030811918678439
$ws.Range("K12").Value = 1.032007189550397
$ws.Range("L12").Value = 1.028224781420289
$ws.Range("M12").Value = 1.036220811077865

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024161238806416
$ws.Range("D13").Value = 1.028427883345574
$ws.Range("E13").Value = 1.024641180259094
$ws.Range("F13").Value = 1.032685171129839
$ws.Range("I13").Value = 1.031793392967034
$ws.Range("J13").Value = 1.030876169909633
$ws.Range("K13").Value = 1.032058920195221
$ws.Range("L13").Value = 1.028286694872374
$ws.Range("M13").Value = 1.036300071516171

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024434270448639
$ws.Range("D14").Value = 1.028629935239422
$ws.Range("E14").Value = 1.024876875946869
$ws.Range("F14").Value = 1.032977607390536
$ws.Range("I14").Value = 1.031869391837529
$ws.Range("J14").Value = 1.031086605972583
$ws.Range("K14").Value = 1.032228329068837
$ws.Range("L14").Value = 1.028489477819897
$ws.Range("M14").Value = 1.036559685054504

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024602432785583
$ws.Range("D15").Value = 1.028754370492155
$ws.Range("E15").Value = 1.025022046144829
$ws.Range("F15").Value = 1.033157729166764
$ws.Range("I15").Value = 1.031916150797161
$ws.Range("J15").Value = 1.031216195496965
$ws.Range("K15").Value = 1.032332638394575
$ws.Range("L15").Value = 1.028614357196013
$ws.Range("M15").Value = 1.036719573493738

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025580281065294
$ws.Range("D16").Value = 1.029477794726654
$ws.Range("E16").Value = 1.025866247957558
$ws.Range("F16").Value = 1.034205242417566
$ws.Range("I16").Value = 1.032187296488926
$ws.Range("J16").Value = 1.031969445428643
$ws.Range("K16").Value = 1.032938714776951
$ws.Range("L16").Value = 1.029340270869527
$ws.Range("M16").Value = 1.037649163692839

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026192851174148
$ws.Range("D17").Value = 1.029930843939386
$ws.Range("E17").Value = 1.026395139460786
$ws.Range("F17").Value = 1.034861561473141
$ws.Range("I17").Value = 1.032356491232462
$ws.Range("J17").Value = 1.032441051077668
$ws.Range("K17").Value = 1.033317972732018
$ws.Range("L17").Value = 1.029794798228848
$ws.Range("M17").Value = 1.038231374077704

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02654986564022
$ws.Range("D18").Value = 1.030194837313329
$ws.Range("E18").Value = 1.026703400337534
$ws.Range("F18").Value = 1.035244112254767
$ws.Range("I18").Value = 1.032454860992907
$ws.Range("J18").Value = 1.032715813628345
$ws.Range("K18").Value = 1.033538859153372
$ws.Range("L18").Value = 1.030059624013456
$ws.Range("M18").Value = 1.03857064815988

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026671549986124
$ws.Range("D19").Value = 1.030284807935947
$ws.Range("E19").Value = 1.026808470169363
$ws.Range("F19").Value = 1.03537450695501
$ws.Range("I19").Value = 1.032488348606146
$ws.Range("J19").Value = 1.032809447064336
$ws.Range("K19").Value = 1.033614120250841
$ws.Range("L19").Value = 1.030149873445198
$ws.Range("M19").Value = 1.038686278043967

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026127158006584
$ws.Range("D20").Value = 1.029882263262796
$ws.Range("E20").Value = 1.026338418535402
$ws.Range("F20").Value = 1.034791172571873
$ws.Range("I20").Value = 1.032338371228029
$ws.Range("J20").Value = 1.032390485116323
$ws.Range("K20").Value = 1.033277315957499
$ws.Range("L20").Value = 1.029746062046677
$ws.Range("M20").Value = 1.038168941536094

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024353882263016
$ws.Range("D21").Value = 1.028570447587395
$ws.Range("E21").Value = 1.024807479811336
$ws.Range("F21").Value = 1.032891504277321
$ws.Range("I21").Value = 1.031847025915259
$ws.Range("J21").Value = 1.031024651788277
$ws.Range("K21").Value = 1.032178456781019
$ws.Range("L21").Value = 1.028429776205474
$ws.Range("M21").Value = 1.036483249532288

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023236279578773
$ws.Range("D22").Value = 1.027743236398251
$ws.Range("E22").Value = 1.023842755926359
$ws.Range("F22").Value = 1.031694592091068
$ws.Range("I22").Value = 1.031535199345164
$ws.Range("J22").Value = 1.0301629769571
$ws.Range("K22").Value = 1.031484550629257
$ws.Range("L22").Value = 1.027599479985435
$ws.Range("M22").Value = 1.035420427492058

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023828999520372
$ws.Range("D23").Value = 1.028181988835758
$ws.Range("E23").Value = 1.024354382482256
$ws.Range("F23").Value = 1.03232934078318
$ws.Range("I23").Value = 1.031700780567697
$ws.Range("J23").Value = 1.030620047151396
$ws.Range("K23").Value = 1.031852691230869
$ws.Range("L23").Value = 1.028039894177053
$ws.Range("M23").Value = 1.035984134250223

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026156842809101
$ws.Range("D24").Value = 1.029904215590868
$ws.Range("E24").Value = 1.026364048998855
$ws.Range("F24").Value = 1.034822979117782
$ws.Range("I24").Value = 1.03234655986279
$ws.Range("J24").Value = 1.032413334677393
$ws.Range("K24").Value = 1.033295688018453
$ws.Range("L24").Value = 1.029768084732949
$ws.Range("M24").Value = 1.038197153101182

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028846947426361
$ws.Range("D25").Value = 1.031892491426988
$ws.Range("E25").Value = 1.028687064456735
$ws.Range("F25").Value = 1.037706200489214
$ws.Range("I25").Value = 1.033083463583557
$ws.Range("J25").Value = 1.034481940473025
$ws.Range("K25").Value = 1.034957351748993
$ws.Range("L25").Value = 1.031762116754378
$ws.Range("M25").Value = 1.040752751152218
